$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace header text and the single data value, then drop the now-unused
# columns (B:D) so only column A ("CT_TX_NEW - Breastfeeding" / 2) remains.
$ws.Range("A1").Value = "CT_TX_NEW - Breastfeeding"
$ws.Range("A2").Value = 2
$ws.Range("B:D").Delete()

# Column A widens to fit the new (longer) header text.
$ws.Columns("A").ColumnWidth = 21.833333333333336

# Match the saved selection/active cell from the workbook.
$ws.Range("B5").Select() | Out-Null
